$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 894954.6
$ws.Range("I15").Value = 894954.6
$ws.Range("K15").Value = 2684863.8
$ws.Range("M15").Value = -2684694.8
$ws.Range("H62").Value = 3060.125
$ws.Range("I62").Value = 2698.6667
$ws.Range("K62").Value = 2698.6667
$ws.Range("M62").Value = -2074.6667
$ws.Range("H65").Value = 3060.125
$ws.Range("I65").Value = 2698.6667
$ws.Range("K65").Value = 13493.3335
$ws.Range("M65").Value = -10373.3335
$ws.Range("H98").Value = 2880.5
$ws.Range("I98").Value = 1296.8889
$ws.Range("J98").Value = 5731
$ws.Range("K98").Value = 1296.8889
$ws.Range("L98").Value = 5731
$ws.Range("M98").Value = 201.1111000000001
$ws.Range("N98").Value = -8727
$ws.Range("H113").Value = 2844.6
$ws.Range("J113").Value = 3222.6667
$ws.Range("L113").Value = 3222.6667
$ws.Range("N113").Value = -9730.6667
$ws.Range("H122").Value = 2880.5
$ws.Range("I122").Value = 1296.8889
$ws.Range("J122").Value = 5731
$ws.Range("K122").Value = 3890.6667
$ws.Range("L122").Value = 17193
$ws.Range("M122").Value = -1440.6667
$ws.Range("N122").Value = -22093
$ws.Range("H132").Value = 1348.8334
$ws.Range("I132").Value = 1348.8334
$ws.Range("K132").Value = 4046.5002
$ws.Range("M132").Value = -1516.5002
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 4707.4043
$ws.Range("I137").Value = 3044.6206
$ws.Range("J137").Value = 7386.3335
$ws.Range("K137").Value = 9133.861800000001
$ws.Range("L137").Value = 22159.0005
$ws.Range("M137").Value = -6583.861800000001
$ws.Range("N137").Value = -27259.0005
$ws.Range("H138").Value = 2681.762
$ws.Range("I138").Value = 2681.762
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8045.286
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -2905.286
$ws.Range("N138").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2755.4521
$ws.Range("I32").Value = 2199.0598
$ws.Range("J32").Value = 8968.5
$ws.Range("K32").Value = 2199.0598
$ws.Range("L32").Value = 8968.5
$ws.Range("M32").Value = -1912.0598
$ws.Range("N32").Value = -9542.5
$ws.Range("H40").Value = 22659.8
$ws.Range("I40").Value = 20000
$ws.Range("K40").Value = 20000
$ws.Range("M40").Value = -19824
$ws.Range("H56").Value = 35000
$ws.Range("J56").Value = 35000
$ws.Range("L56").Value = 35000
$ws.Range("N56").Value = -36484
$ws.Range("H102").Value = 2200.348
$ws.Range("I102").Value = 1845.8182
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1845.8182
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -223.8181999999999
$ws.Range("N102").Value = -13244
$ws.Range("H122").Value = 4601.0835
$ws.Range("I122").Value = 3553.3333
$ws.Range("J122").Value = 5648.8335
$ws.Range("K122").Value = 10659.9999
$ws.Range("L122").Value = 16946.5005
$ws.Range("M122").Value = -8209.999899999999
$ws.Range("N122").Value = -21846.5005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 6534.5
$ws.Range("J25").Value = 12069
$ws.Range("L25").Value = 12069
$ws.Range("N25").Value = -12539
$ws.Range("H76").Value = 178136.17
$ws.Range("J76").Value = 178136.17
$ws.Range("L76").Value = 178136.17
$ws.Range("N76").Value = -178766.17
$ws.Range("H79").Value = 178136.17
$ws.Range("J79").Value = 178136.17
$ws.Range("L79").Value = 178136.17
$ws.Range("N79").Value = -180320.17
$ws.Range("H134").Value = 3253.8
$ws.Range("I134").Value = 1615.3334
$ws.Range("K134").Value = 4846.0002
$ws.Range("M134").Value = -2311.0002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 995
$ws.Range("J29").Value = 995
$ws.Range("L29").Value = 995
$ws.Range("N29").Value = -1581
$ws.Range("H43").Value = 146034.14
$ws.Range("J43").Value = 146034.14
$ws.Range("L43").Value = 146034.14
$ws.Range("N43").Value = -146402.14
$ws.Range("H58").Value = 5318.2104
$ws.Range("I58").Value = 2780.6155
$ws.Range("K58").Value = 2780.6155
$ws.Range("M58").Value = -2577.6155
$ws.Range("H86").Value = 6858.5713
$ws.Range("I86").Value = 6002
$ws.Range("K86").Value = 6002
$ws.Range("M86").Value = -4879
$ws.Range("H89").Value = 6858.5713
$ws.Range("I89").Value = 6002
$ws.Range("K89").Value = 30010
$ws.Range("M89").Value = -24394
$ws.Range("H94").Value = 2465.4666
$ws.Range("I94").Value = 1456.375
$ws.Range("K94").Value = 1456.375
$ws.Range("M94").Value = -1005.375
$ws.Range("H101").Value = 146034.14
$ws.Range("J101").Value = 146034.14
$ws.Range("L101").Value = 146034.14
$ws.Range("N101").Value = -152524.14
$ws.Range("H105").Value = 2313.037
$ws.Range("I105").Value = 1978.85
$ws.Range("K105").Value = 1978.85
$ws.Range("M105").Value = -231.8499999999999
$ws.Range("H132").Value = 2765.205
$ws.Range("I132").Value = 2364.2896
$ws.Range("K132").Value = 7092.8688
$ws.Range("M132").Value = -4562.8688
$ws.Range("H136").Value = 5318.2104
$ws.Range("I136").Value = 2780.6155
$ws.Range("K136").Value = 8341.8465
$ws.Range("M136").Value = -5791.8465

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 399
$ws.Range("I8").Value = 399
$ws.Range("K8").Value = 1197
$ws.Range("M8").Value = -1058
$ws.Range("H97").Value = 4058
$ws.Range("J97").Value = 4549.8184
$ws.Range("L97").Value = 13649.4552
$ws.Range("N97").Value = -14641.4552
$ws.Range("H141").Value = 4448.5
$ws.Range("I141").Value = 5297.8
$ws.Range("K141").Value = 15893.4
$ws.Range("M141").Value = -10713.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H59").Value = 25000
$ws.Range("I59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("K59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("M59").Value = -24417
$ws.Range("N59").Value = -26166
$ws.Range("H69").Value = 174050
$ws.Range("J69").Value = 174050
$ws.Range("L69").Value = 174050
$ws.Range("N69").Value = -175548
$ws.Range("H72").Value = 174050
$ws.Range("J72").Value = 174050
$ws.Range("L72").Value = 522150
$ws.Range("N72").Value = -529638
$ws.Range("H97").Value = 759.44446
$ws.Range("I97").Value = 624.8
$ws.Range("J97").Value = 1432.6666
$ws.Range("K97").Value = 624.8
$ws.Range("L97").Value = 1432.6666
$ws.Range("M97").Value = -128.8
$ws.Range("N97").Value = -2424.6666
$ws.Range("H107").Value = 482.66666
$ws.Range("I107").Value = 399
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 399
$ws.Range("L107").Value = 650
$ws.Range("M107").Value = 1521
$ws.Range("N107").Value = -4490
$ws.Range("H113").Value = 3383.5518
$ws.Range("I113").Value = 2665.9473
$ws.Range("K113").Value = 2665.9473
$ws.Range("M113").Value = -495.9472999999998
$ws.Range("H132").Value = 3640.111
$ws.Range("I132").Value = 2912.9092
$ws.Range("J132").Value = 6839.8
$ws.Range("K132").Value = 8738.7276
$ws.Range("L132").Value = 20519.4
$ws.Range("M132").Value = -6208.7276
$ws.Range("N132").Value = -25579.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2812.44
$ws.Range("I46").Value = 1523.7858
$ws.Range("J46").Value = 4452.5454
$ws.Range("K46").Value = 1523.7858
$ws.Range("L46").Value = 4452.5454
$ws.Range("M46").Value = -1335.7858
$ws.Range("N46").Value = -4828.5454

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9451
$ws.Range("I62").Value = 9033
$ws.Range("J62").Value = 9764.5
$ws.Range("K62").Value = 9033
$ws.Range("L62").Value = 9764.5
$ws.Range("M62").Value = -8409
$ws.Range("N62").Value = -11012.5
$ws.Range("H65").Value = 9451
$ws.Range("I65").Value = 9033
$ws.Range("J65").Value = 9764.5
$ws.Range("K65").Value = 45165
$ws.Range("L65").Value = 48822.5
$ws.Range("M65").Value = -42045
$ws.Range("N65").Value = -55062.5
$ws.Range("H126").Value = 4223.75
$ws.Range("I126").Value = 2966.3333
$ws.Range("J126").Value = 4978.2
$ws.Range("K126").Value = 8898.999899999999
$ws.Range("L126").Value = 14934.6
$ws.Range("M126").Value = -6428.999899999999
$ws.Range("N126").Value = -19874.6
$ws.Range("H132").Value = 5426.8
$ws.Range("I132").Value = 3921.4092
$ws.Range("J132").Value = 7974.385
$ws.Range("K132").Value = 11764.2276
$ws.Range("L132").Value = 23923.155
$ws.Range("M132").Value = -9234.2276
$ws.Range("N132").Value = -28983.155
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
